$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Old:  A1=項目  B1=日期        C1=提前提醒
# New:  A1=項目名稱 B1=日期      C1=提前提醒  (A1 text changed)
$ws.Range("A1").Value = "項目名稱"
$ws.Range("B1").Value = "日期"
$ws.Range("C1").Value = "提前提醒"

# --- Data row (row 2) ---------------------------------------------------
# Old: A2=結婚紀念日 B2=2022/3/18 (44638) C2=7
# New: A2=結婚紀念日 B2=2022/5/11 (44692) C2=2
$ws.Range("A2").Value = "結婚紀念日"
$ws.Range("B2").Value = 44692
$ws.Range("C2").Value = 2

# --- Remove the old third row (xxx / 44690 / 1) --------------------------
$ws.Rows("3").Delete()

# --- Number format for the date column -----------------------------------
$ws.Range("B1:B2").NumberFormat = "mm/dd;@"

# --- Font change for the whole used range ---------------------------------
$ws.Range("A1:C2").Font.Name = "Adobe 繁黑體 Std B"

# --- Page setup (paper size + orientation) --------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection -------------------------------------------------------------
[void]$ws.Range("G15").Select()

Write-Host "done"
